$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates that Excel would otherwise auto-convert to numbers;
# force text format first so the literal string is preserved, matching the
# source data which stores these as inline strings.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.47"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2637"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06360"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.50"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07834"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.566"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5522"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.60"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.653"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.77"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.19"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.020"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.89"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1221"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.180"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05885"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.279"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.577"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.268"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.611"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.823"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9585"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.426"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5773"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01600"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8623"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.833"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.008"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.967"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05164"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.429"

# Price column (D) updates that are already non-numeric-looking text (contain
# multiple separators or special glyphs), so plain assignment keeps them as text.
$ws.Range("D2").Value = "26.273.41"
$ws.Range("D3").Value = "1.662.44"
$ws.Range("D13").Value = "1.662.51"
$ws.Range("D14").Value = "1.890.04"
$ws.Range("D16").Value = "0.0₅8158"
$ws.Range("D42").Value = "1.042.60"
$ws.Range("D44").Value = "1.802.62"

# Volume(1h) column (E) updates - always padded percentage text, safe to assign directly.
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  +0.78%  "
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("E24").Value = "  +2.49%  "
$ws.Range("E25").Value = "  -2.39%  "
$ws.Range("E26").Value = "  -1.50%  "
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("E28").Value = "  +2.88%  "
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  +1.43%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("E39").Value = "  +1.74%  "
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E46").Value = "  -4.83%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("E51").Value = "  -4.25%  "
